# Updates cryptos list values (Price and Volume(1h) columns) to match
# the latest scrape, per commit "Updated cryptos list on Mon Oct 16
# 20:40:37 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.546.07'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.77%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.592.56'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.44%  '

$ws.Range("E4").Value = '  -0.42%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.78%  '

$ws.Range("E6").Value = '  +1.30%  '

$ws.Range("E7").Value = '  -0.52%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.03'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +8.91%  '

$ws.Range("E9").Value = '  +1.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0602'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.62%  '

$ws.Range("E11").Value = '  +2.30%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.819.96'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.42%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.584.57'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.14%  '

$ws.Range("E14").Value = '  +0.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.531'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.29%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.533.18'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.89%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.04'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.94%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '232.62'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.58%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.52'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.60%  '

$ws.Range("E20").Value = '  +1.14%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.55%  '

$ws.Range("E22").Value = '  -0.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.43'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.37%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.14%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.84'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.28%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.33'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.71%  '

$ws.Range("E27").Value = '  +0.06%  '

$ws.Range("E28").Value = '  +1.08%  '

$ws.Range("E29").Value = '  -0.50%  '

$ws.Range("E30").Value = '  +0.70%  '

$ws.Range("E31").Value = '  +0.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.25'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.48%  '

$ws.Range("E33").Value = '  -0.47%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.416.74'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.19%  '

$ws.Range("E35").Value = '  -0.98%  '

$ws.Range("E36").Value = '  -6.36%  '

$ws.Range("E37").Value = '  -0.21%  '

$ws.Range("E38").Value = '  +0.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.55'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.74%  '

$ws.Range("E40").Value = '  +2.21%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.817'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.80%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.74'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.56%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.980'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.11%  '

$ws.Range("E45").Value = '  +6.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.75'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.733.09'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.65%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.96'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.18%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0106'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.06%  '

$ws.Range("E50").Value = '  -0.17%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '39.88'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +17.56%  '
